# Apply updated plate-fit coefficients / recomputed variance columns
# (values taken from the recalculated least-squares fit; the "variance_*"
# columns G-J shrink by several orders of magnitude because the new code
# computes them correctly, while columns C-F keep the same coefficients
# modulo floating-point solver noise from the recompute.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
  @("C2", "-7.409452811607048e-15"),
  @("D2", "5.179118952131721e-14"),
  @("E2", "-1.293193473336941e-16"),
  @("F2", "9.039267806713712e-16"),
  @("G2", "1.17828454229957e-28"),
  @("H2", "3.589259970498579e-32"),
  @("I2", "2.88157300524464e-28"),
  @("J2", "8.777773337846538e-32"),
  @("C3", "-4.861914998617235e-14"),
  @("D3", "-4.881752797759572e-15"),
  @("E3", "-8.485642467796631e-16"),
  @("F3", "-8.520265958934939e-17"),
  @("G3", "1.039796723505734e-28"),
  @("H3", "3.167401950169903e-32"),
  @("I3", "2.542891858318357e-28"),
  @("J3", "7.746091566775665e-32"),
  @("C4", "2.243882269273335e-14"),
  @("D4", "8.207440763935868e-14"),
  @("E4", "3.916313362594168e-16"),
  @("F4", "1.432468644930796e-15"),
  @("G4", "1.873639817155104e-28"),
  @("H4", "5.707433267114283e-32"),
  @("I4", "4.582110453667471e-28"),
  @("J4", "1.39579066357389e-31"),
  @("C5", "-4.755730566440148e-14"),
  @("D5", "-3.604483617815153e-14"),
  @("E5", "-8.300315672211551e-16"),
  @("F5", "-6.291010696507135e-16"),
  @("G5", "1.094360473638643e-27"),
  @("H5", "3.333612637963637e-31"),
  @("I5", "2.676331128548492e-27"),
  @("J5", "8.152570828732916e-31"),
  @("C6", "3.046466708966163e-11"),
  @("D6", "4.759676939831958e-11"),
  @("E6", "5.317087462385539e-13"),
  @("F6", "8.30720339313157e-13"),
  @("G6", "1.438220869630538e-22"),
  @("H6", "4.381071303902507e-26"),
  @("I6", "3.517264535626203e-22"),
  @("J6", "1.07142004754975e-25"),
  @("C7", "-3.94723535516304e-10"),
  @("D7", "-6.836102541605435e-11"),
  @("E7", "-6.88922532987228e-12"),
  @("F7", "-1.193124973549675e-12"),
  @("G7", "2.55345213899749e-22"),
  @("H7", "7.778260021302673e-26"),
  @("I7", "6.244636579513614e-22"),
  @("J7", "1.902225082357134e-25"),
  @("C8", "1.58365436615912e-11"),
  @("D8", "-2.078024447514249e-10"),
  @("E8", "2.763998290306051e-13"),
  @("F8", "-3.626836854605974e-12"),
  @("G8", "4.636377349421205e-22"),
  @("H8", "1.412321305338226e-25"),
  @("I8", "1.133856834457491e-21"),
  @("J8", "3.45392543319966e-25"),
  @("C9", "1.090878687667711e-05"),
  @("D9", "-0.0005432880894939687"),
  @("E9", "1.903942483963642e-07"),
  @("F9", "-9.482165948539368e-06"),
  @("G9", "1.926360197650284e-15"),
  @("H9", "5.868028729880435e-19"),
  @("I9", "4.711041641175583e-15"),
  @("J9", "1.435065349242647e-18"),
  @("C10", "0.0005447916065334166"),
  @("D10", "1.09788290723787e-05"),
  @("E10", "9.508407271237571e-06"),
  @("F10", "1.916167153266832e-07"),
  @("G10", "5.817463482871495e-15"),
  @("H10", "1.772100715855714e-18"),
  @("I10", "1.422699282681159e-14"),
  @("J10", "4.333789846227358e-18"),
  @("D11", "4.392311385722874"),
  @("G11", "1.534528052726398e-09"),
  @("H11", "4.674439760118376e-13"),
  @("I11", "3.752790140060064e-09"),
  @("J11", "1.143165249466097e-12")
)

foreach ($u in $updates) {
  $ws.Range($u[0]).Value = [double]$u[1]
}

Write-Output "updated $($updates.Count) cells"